$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Semilla 2")

# --- Cell value corrections (rows 9-14) -------------------------------
# Row 9
$ws.Range("C9").Value = "3045984642"
$ws.Range("E9").Value = "3046010523"

# Row 10
$ws.Range("E10").Value = "3045981670"

# Row 11
$ws.Range("C11").Value = "3043118820"
$ws.Range("D11").Value = "732111324706350"

# Row 12
$ws.Range("C12").Value = "3045984556"

# Row 14
$ws.Range("B14").Value = "582710820"
$ws.Range("C14").Value = "3045984642"
$ws.Range("D14").Value = "732111324707276"

# --- Formatting touch-ups ----------------------------------------------
# E11:I12 lose their inherited number format / pick up the underlined font
$ws.Range("E11:I12").Font.Underline = $true
$ws.Range("E13:I14").Font.Underline = $true

# --- Sheet view ---------------------------------------------------------
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Range("D12").Select()
